$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24
$ws.Range("D24").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E24").Value = "['HardwareFault']"

# Row 26
$ws.Range("D26").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['SoftwareFault']"

# Row 27
$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

# Row 56
$ws.Range("D56").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "[]"

# Row 69
$ws.Range("D69").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E69").Value = "['Normal', 'SurroundingEnvironment']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'HardwareFault']"

# Row 74
$ws.Range("D74").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E74").Value = "['Normal', 'SoftwareFault']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E113").Value = "['Normal', 'SoftwareFault']"
